$wb = $excel.ActiveWorkbook

# --- Original sheet "7.0-9.3" becomes "7.0-7.5" (keeps sheetId=6) ---
$ws6 = $wb.Worksheets.Item("7.0-9.3")

# Create the three new tabs by copying the existing sheet (preserves header
# style / formatting) and placing each right after the previous one, so the
# final tab order is: ... , 7.0-7.5, 7.6-7.11, 8.0-8.5, 8.6-9.3
$ws6.Copy($null, $ws6)
$ws7 = $wb.Worksheets.Item("7.0-9.3 (2)")
$ws7.Name = "7.6-7.11"

$ws7.Copy($null, $ws7)
$ws8 = $wb.Worksheets.Item("7.6-7.11 (2)")
$ws8.Name = "8.0-8.5"

$ws8.Copy($null, $ws8)
$ws9 = $wb.Worksheets.Item("8.0-8.5 (2)")
$ws9.Name = "8.6-9.3"

# Rename the original sheet last (so the lookups above by old name still work)
$ws6.Name = "7.0-7.5"

# --- Updated "ss" (column B) values for each tab ---

$vals6 = @(47,49,51,53,55,57,59,61,62,64,66,68,69,71,73,75,76,78,80,82,84,85,87,89,91,93,94,96,98,100,102,104,106,109,111,113,116,118)
$vals7 = @(43,45,47,48,50,52,54,55,57,59,60,62,63,65,67,68,70,72,73,75,77,78,80,82,84,85,87,89,91,93,95,97,99,102,104,107,110,113)
$vals8 = @(40,41,43,44,46,48,49,51,52,54,55,57,58,60,61,63,65,66,68,69,71,73,74,76,78,80,81,83,85,87,89,92,94,97,100,104,118,118)
$vals9 = @(40,40,40,40,41,43,44,46,47,48,50,51,53,54,56,57,59,60,62,63,65,66,68,70,71,73,75,77,79,113,113,113,113,113,113,113,113,113)

for ($i = 0; $i -lt $vals6.Length; $i++) {
    $ws6.Cells.Item($i + 2, 2).Value = $vals6[$i]
}
for ($i = 0; $i -lt $vals7.Length; $i++) {
    $ws7.Cells.Item($i + 2, 2).Value = $vals7[$i]
}
for ($i = 0; $i -lt $vals8.Length; $i++) {
    $ws8.Cells.Item($i + 2, 2).Value = $vals8[$i]
}
for ($i = 0; $i -lt $vals9.Length; $i++) {
    $ws9.Cells.Item($i + 2, 2).Value = $vals9[$i]
}

# Restore the originally-active tab (first sheet), since our copy/rename
# operations above shift Excel's "active sheet" cursor as a side effect.
$wb.Worksheets.Item("5.0-5.3").Activate()
